$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) stays the same text, just re-write to keep things
# explicit / normalised after the shared-string table is rebuilt) ----
$ws.Range("A1").Value = "sr"
$ws.Range("B1").Value = "github_username"
$ws.Range("C1").Value = "repo_name_to_import"
$ws.Range("D1").Value = "azure_target_namespace"

# ---- New repo list (9 rows) ----
$repoNames = @(
    "casaplotms ",
    "casa-build-utils ",
    "casashell",
    "casaaddons ",
    "cartavis ",
    "carta-casacore",
    "casa-asap ",
    "almatasks ",
    "app-n-pak "
)

# Column A - serial numbers, written top to bottom
for ($i = 0; $i -lt $repoNames.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i + 1
}

# Column B - github username, same value filled down the whole column
$ws.Range("B2:B10").Value = "code-migration"

# Column C - repo name to import, one distinct value per row
for ($i = 0; $i -lt $repoNames.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $repoNames[$i]
}

# Column D - azure target namespace, same value filled down the whole column
$ws.Range("D2:D10").Value = "repo-migartion/git-project"

# ---- C4 ("casashell") picked up hyperlink-like formatting (underline
# removed, hyperlink theme colour, left/center aligned, wrapped text) ----
$c4 = $ws.Range("C4")
$c4.Style = "Hyperlink"
$c4.Font.Underline = $false
$c4.Font.ThemeColor = 11
$c4.WrapText = $true
$c4.HorizontalAlignment = -4131
$c4.VerticalAlignment = -4108

# ---- Selection moves to C12 ----
$ws.Range("C12").Select() | Out-Null
